# Add a new CV "award" entry: Data Science Fellow (2021).
# This inserts a new row at row 32 (pushing the existing rows below it down
# by one) and populates it with the new award entry, matching the author's
# "add data science fellow" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 32 ("Outstanding reviewer...").
$ws.Rows("32:32").Insert()

# Populate the new row 32 with the Data Science Fellow award entry.
$ws.Range("A32").Value2 = "award"
$ws.Range("C32").Value2 = 2021
$ws.Range("D32").Value2 = "*Data Science Fellow*: College of Education, University of Oregon"

# The wrapped "what" text needs the taller (two-line) row height, matching
# the other multi-line entries elsewhere in the sheet.
$ws.Rows("32:32").RowHeight = 34

# Match the workbook's final selection state (D32 selected).
[void]$ws.Range("D32").Select()
